$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 2378
$ws.Range("E2").Value = 69
$ws.Range("F2").Value = 69
$ws.Range("G2").Value = 47
$ws.Range("H2").Value = 31
$ws.Range("I2").Value = 28
$ws.Range("J2").Value = 3
$ws.Range("K2").Value = 2453
$ws.Range("L2").Value = 1184
$ws.Range("M2").Value = 1269
$ws.Range("N2").Value = 1293
$ws.Range("O2").Value = -25
$ws.Range("P2").Value = 626
$ws.Range("Q2").Value = 169
$ws.Range("R2").Value = -77
$ws.Range("S2").Value = -97
$ws.Range("T2").Value = 37
$ws.Range("U2").Value = 132
$ws.Range("V2").Value = 893
$ws.Range("W2").Value = 2.88
$ws.Range("X2").Value = 1.3
$ws.Range("Y2").Value = 2.17
$ws.Range("Z2").Value = 1.27
$ws.Range("AA2").Value = 93.31
$ws.Range("AB2").Value = 118.49
$ws.Range("AC2").Value = 48
$ws.Range("AD2").Value = 26.07
$ws.Range("AE2").Value = 2564
$ws.Range("AF2").Value = 0.49
$ws.Range("AG2").Value = 20
$ws.Range("AH2").Value = 1.59
$ws.Range("AI2").Value = 36.61
$ws.Range("AJ2").Value = 57000000

# Row 3
$ws.Range("D3").Value = 2917
$ws.Range("E3").Value = 98
$ws.Range("F3").Value = 98
$ws.Range("G3").Value = 44
$ws.Range("H3").Value = 42
$ws.Range("I3").Value = 40
$ws.Range("J3").Value = 2
$ws.Range("K3").Value = 3050
$ws.Range("L3").Value = 1708
$ws.Range("M3").Value = 1342
$ws.Range("N3").Value = 1365
$ws.Range("O3").Value = -23
$ws.Range("P3").Value = 626
$ws.Range("Q3").Value = 17
$ws.Range("R3").Value = -222
$ws.Range("S3").Value = 228
$ws.Range("T3").Value = 28
$ws.Range("U3").Value = -11
$ws.Range("V3").Value = 1168
$ws.Range("W3").Value = 3.37
$ws.Range("X3").Value = 1.44
$ws.Range("Y3").Value = 3.03
$ws.Range("Z3").Value = 1.52
$ws.Range("AA3").Value = 127.27
$ws.Range("AB3").Value = 123.28
$ws.Range("AC3").Value = 71
$ws.Range("AD3").Value = 30.11
$ws.Range("AE3").Value = 2706
$ws.Range("AF3").Value = 0.79
$ws.Range("AG3").Value = 20
$ws.Range("AH3").Value = 0.94
$ws.Range("AI3").Value = 25.02
$ws.Range("AJ3").Value = 57000000

# Row 4
$ws.Range("D4").Value = 3108
$ws.Range("E4").Value = 101
$ws.Range("F4").Value = 101
$ws.Range("G4").Value = 53
$ws.Range("H4").Value = 40
$ws.Range("I4").Value = 39
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 3415
$ws.Range("L4").Value = 2023
$ws.Range("M4").Value = 1393
$ws.Range("N4").Value = 1415
$ws.Range("O4").Value = -22
$ws.Range("P4").Value = 626
$ws.Range("Q4").Value = 44
$ws.Range("R4").Value = -189
$ws.Range("S4").Value = 135
$ws.Range("T4").Value = 77
$ws.Range("U4").Value = -33
$ws.Range("V4").Value = 1346
$ws.Range("W4").Value = 3.26
$ws.Range("X4").Value = 1.28
$ws.Range("Y4").Value = 2.8
$ws.Range("Z4").Value = 1.23
$ws.Range("AA4").Value = 145.27
$ws.Range("AB4").Value = 127.83
$ws.Range("AC4").Value = 68
$ws.Range("AD4").Value = 46.6
$ws.Range("AE4").Value = 2805
$ws.Range("AF4").Value = 1.14
$ws.Range("AG4").Value = 0
$ws.Range("AH4").Value = 0
$ws.Range("AI4").Value = 0
$ws.Range("AJ4").Value = 57000000

# Row 5
$ws.Range("D5").Value = 2910
$ws.Range("E5").Value = 70
$ws.Range("F5").Value = 70
$ws.Range("G5").Value = 17
$ws.Range("H5").Value = 7
$ws.Range("I5").Value = 6
$ws.Range("J5").Value = 1
$ws.Range("K5").Value = 3079
$ws.Range("L5").Value = 1761
$ws.Range("M5").Value = 1318
$ws.Range("N5").Value = 1339
$ws.Range("O5").Value = -21
$ws.Range("P5").Value = 626
$ws.Range("Q5").Value = 111
$ws.Range("R5").Value = -108
$ws.Range("S5").Value = 20
$ws.Range("T5").Value = 38
$ws.Range("U5").Value = 73
$ws.Range("V5").Value = 1289
$ws.Range("W5").Value = 2.41
$ws.Range("X5").Value = 0.22
$ws.Range("Y5").Value = 0.42
$ws.Range("Z5").Value = 0.2
$ws.Range("AA5").Value = 133.63
$ws.Range("AB5").Value = 129.28
$ws.Range("AC5").Value = 10
$ws.Range("AD5").Value = 159.6
$ws.Range("AE5").Value = 2656
$ws.Range("AF5").Value = 0.61
$ws.Range("AG5").Value = 0
$ws.Range("AH5").Value = 0
$ws.Range("AI5").Value = 0
$ws.Range("AJ5").Value = 57000000

# Row 6
$ws.Range("D6").Value = 2777
$ws.Range("E6").Value = 34
$ws.Range("F6").Value = 34
$ws.Range("G6").Value = -41
$ws.Range("H6").Value = -66
$ws.Range("I6").Value = -67
$ws.Range("K6").Value = 3048
$ws.Range("L6").Value = 1728
$ws.Range("M6").Value = 1321
$ws.Range("N6").Value = 1342
$ws.Range("P6").Value = 626
$ws.Range("Q6").Value = -15
$ws.Range("R6").Value = -42
$ws.Range("S6").Value = 38
$ws.Range("T6").Value = 39
$ws.Range("U6").Value = -54
$ws.Range("V6").Value = 1298
$ws.Range("W6").Value = 1.23
$ws.Range("X6").Value = -2.38
$ws.Range("Y6").Value = -4.97
$ws.Range("Z6").Value = -2.16
$ws.Range("AA6").Value = 130.79
$ws.Range("AB6").Value = 123.85
$ws.Range("AC6").Value = -117
$ws.Range("AD6").Value = -11.24
$ws.Range("AE6").Value = 2660
$ws.Range("AF6").Value = 0.49
$ws.Range("AG6").Value = 0
$ws.Range("AH6").Value = 0
$ws.Range("AI6").Value = 0
$ws.Range("AJ6").Value = 57000000

# Clear rows 7-9 columns D:AI (keep A,B,C)
$ws.Range("D7:AI7").ClearContents()
$ws.Range("D8:AI8").ClearContents()
$ws.Range("D9:AI9").ClearContents()